$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9082
$ws.Range("F3").Value = 1982
$ws.Range("F4").Value = 6639
$ws.Range("F5").Value = 180
$ws.Range("F6").Value = 2156
$ws.Range("F7").Value = 609
$ws.Range("F10").Value = 77
$ws.Range("F13").Value = 14
$ws.Range("F14").Value = 86
$ws.Range("F15").Value = 29
$ws.Range("F16").Value = 9018
$ws.Range("F21").Value = 1850
$ws.Range("F25").Value = 105
$ws.Range("F27").Value = 204
$ws.Range("F28").Value = 1048
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 78
$ws.Range("F31").Value = 567
$ws.Range("F32").Value = 37
$ws.Range("F33").Value = 47
$ws.Range("F35").Value = 2394
$ws.Range("F36").Value = 885
$ws.Range("F37").Value = 559
$ws.Range("F41").Value = 312
$ws.Range("F42").Value = 186
$ws.Range("F44").Value = 100
$ws.Range("F45").Value = 32
$ws.Range("F46").Value = 87
$ws.Range("F47").Value = 23
$ws.Range("F48").Value = 4005
$ws.Range("F49").Value = 13

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2344
$ws.Range("F4").Value = 347
$ws.Range("F5").Value = 27

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2344
$ws.Range("F3").Value = 9082
$ws.Range("F5").Value = 1982
$ws.Range("F6").Value = 6639
$ws.Range("F8").Value = 609
$ws.Range("F13").Value = 77
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 14
$ws.Range("F16").Value = 86
$ws.Range("F17").Value = 9018
$ws.Range("F22").Value = 1850
$ws.Range("F24").Value = 105
$ws.Range("F26").Value = 204
$ws.Range("F27").Value = 20
$ws.Range("F29").Value = 567
$ws.Range("F30").Value = 37
$ws.Range("F31").Value = 47
$ws.Range("F33").Value = 885
$ws.Range("F34").Value = 21
$ws.Range("F36").Value = 559
$ws.Range("F37").Value = 312
$ws.Range("F39").Value = 186
$ws.Range("F41").Value = 100
$ws.Range("F42").Value = 32
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 23
$ws.Range("F45").Value = 4005
$ws.Range("F48").Value = 13
